$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - copy style from E1 (bold header style) then set value
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Data cells F2:F30 - plain values, no special style
$ws.Range("F2").Value = "2021-10-05 10:50:57.235376"
$ws.Range("F3").Value = "2021-10-05 10:50:57.235389"
$ws.Range("F4").Value = "2021-10-05 10:50:57.235393"
$ws.Range("F5").Value = "2021-10-05 10:50:57.235396"
$ws.Range("F6").Value = "2021-10-05 10:50:57.235400"
$ws.Range("F7").Value = "2021-10-05 10:50:57.235403"
$ws.Range("F8").Value = "2021-10-05 10:50:57.235406"
$ws.Range("F9").Value = "2021-10-05 10:50:57.235409"
$ws.Range("F10").Value = "2021-10-05 10:50:57.235413"
$ws.Range("F11").Value = "2021-10-05 10:50:57.235416"
$ws.Range("F12").Value = "2021-10-05 10:50:57.235419"
$ws.Range("F13").Value = "2021-10-05 10:50:57.235422"
$ws.Range("F14").Value = "2021-10-05 10:50:57.235425"
$ws.Range("F15").Value = "2021-10-05 10:50:57.235428"
$ws.Range("F16").Value = "2021-10-05 10:50:57.235431"
$ws.Range("F17").Value = "2021-10-05 10:50:57.235434"
$ws.Range("F18").Value = "2021-10-05 10:50:57.235437"
$ws.Range("F19").Value = "2021-10-05 10:50:57.235440"
$ws.Range("F20").Value = "2021-10-05 10:50:57.235443"
$ws.Range("F21").Value = "2021-10-05 10:50:57.235447"
$ws.Range("F22").Value = "2021-10-05 10:50:57.235450"
$ws.Range("F23").Value = "2021-10-05 10:50:57.235452"
$ws.Range("F24").Value = "2021-10-05 10:50:57.235456"
$ws.Range("F25").Value = "2021-10-05 10:50:57.235459"
$ws.Range("F26").Value = "2021-10-05 10:50:57.235462"
$ws.Range("F27").Value = "2021-10-05 10:50:57.235465"
$ws.Range("F28").Value = "2021-10-05 10:50:57.235468"
$ws.Range("F29").Value = "2021-10-05 10:50:57.235471"
$ws.Range("F30").Value = "2021-10-05 10:50:57.235474"
